$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("191:192").Insert()

$ws.Range("A191").Value = 4
$ws.Range("B191").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C191").Value = "Los Lagos"
$ws.Range("D191").Value = 44488
$ws.Range("E191").Value = 10
$ws.Range("F191").Value = "Fruta"
$ws.Range("G191").Value = 100102
$ws.Range("H191").Value = "Cítricos"
$ws.Range("I191").Value = 100102003
$ws.Range("J191").Value = "Limón"
$ws.Range("K191").Value = "Sin especificar"
$ws.Range("L191").Value = "1a amarillo"
$ws.Range("M191").Value = 1200
$ws.Range("N191").Value = 9000
$ws.Range("O191").Value = 9500
$ws.Range("P191").Value = 9250
$ws.Range("Q191").Value = "$/malla 18 kilos"
$ws.Range("R191").Value = "Provincia de Melipilla"
$ws.Range("S191").Value = 514
$ws.Range("T191").Value = 18

$ws.Range("A192").Value = 4
$ws.Range("B192").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C192").Value = "Los Lagos"
$ws.Range("D192").Value = 44488
$ws.Range("E192").Value = 10
$ws.Range("F192").Value = "Fruta"
$ws.Range("G192").Value = 100102
$ws.Range("H192").Value = "Cítricos"
$ws.Range("I192").Value = 100102003
$ws.Range("J192").Value = "Limón"
$ws.Range("K192").Value = "Sin especificar"
$ws.Range("L192").Value = "2a amarillo"
$ws.Range("M192").Value = 600
$ws.Range("N192").Value = 7500
$ws.Range("O192").Value = 7500
$ws.Range("P192").Value = 7500
$ws.Range("Q192").Value = "$/malla 18 kilos"
$ws.Range("R192").Value = "Provincia de Melipilla"
$ws.Range("S192").Value = 417
$ws.Range("T192").Value = 18
